# "Disaggregation of commodity Copper"
#
# Each yearly worksheet (2000 .. 2100) holds a small 4-row table
# (Neodymium / Dysprosium / Copper / Raw silicon) x 3 sector columns
# (D = Photovoltaic plants, E = Offshore wind plants, F = Onshore wind
# plants), in rows 5 (Neodymium), 6 (Dysprosium), 7 (Copper), 8 (Raw
# silicon).
#
# The commit renames the "Copper ores and concentrates" label to plain
# "Copper" and re-distributes the Copper/Neodymium/Raw-silicon figures
# across the three sector columns: for rows 5, 7 and 8 the D/E/F values
# are rotated one column to the right (D->E, E->F, F->D), i.e.
#   new D = old F
#   new E = old D
#   new F = old E
# Row 6 (Dysprosium) is left untouched.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename the commodity label in column C, row 7, from
    # "Copper ores and concentrates" to "Copper".
    $ws.Range("C7").Value = "Copper"

    # Rotate the D/E/F (sector) values one column to the right for the
    # Neodymium (5), Copper (7) and Raw silicon (8) rows.
    foreach ($row in 5, 7, 8) {
        $dCell = $ws.Range("D$row")
        $eCell = $ws.Range("E$row")
        $fCell = $ws.Range("F$row")

        $dVal = $dCell.Value()
        $eVal = $eCell.Value()
        $fVal = $fCell.Value()

        $dCell.Value = $fVal
        $eCell.Value = $dVal
        $fCell.Value = $eVal
    }
}
